$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 25: "Splash screen" status -> highlight (same fill style as header cell F1, s="8") + progress 0
$ws.Range("F1").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("E25").Value = 0

# Row 27: bump progress to 90% and update comment text
$ws.Range("E27").Value = 0.9
$ws.Range("F27").Value = "Quasi fini : paufinage + retours"

# Row 28: highlight like row 25/header, set progress style + value, add comment
$ws.Range("F1").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0.9
$ws.Range("F28").Value = "Quasi fini : paufinage + retours"

# Row 63: bump progress to 90% and update comment text
$ws.Range("E63").Value = 0.9
$ws.Range("F63").Value = "D'autres trucs à rajouter par la suite ?"

# Move selection to reflect the latest edit location
$ws.Range("B60").Select()
